$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 475.7143
$ws.Range("I8").Value = 262
$ws.Range("K8").Value = 786
$ws.Range("M8").Value = -647
$ws.Range("H15").Value = 396.36
$ws.Range("I15").Value = 396.36
$ws.Range("K15").Value = 1189.08
$ws.Range("M15").Value = -1020.08
$ws.Range("H131").Value = 1910.6154
$ws.Range("I131").Value = 1685.4
$ws.Range("J131").Value = 2661.3333
$ws.Range("K131").Value = 5056.200000000001
$ws.Range("L131").Value = 7983.999899999999
$ws.Range("M131").Value = -16.20000000000073
$ws.Range("N131").Value = -18063.9999
$ws.Range("H137").Value = 1328.5209
$ws.Range("I137").Value = 1071.6389
$ws.Range("J137").Value = 2099.1667
$ws.Range("K137").Value = 3214.9167
$ws.Range("L137").Value = 6297.500100000001
$ws.Range("M137").Value = -664.9166999999998
$ws.Range("N137").Value = -11397.5001
$ws.Range("H138").Value = 2370.3193
$ws.Range("I138").Value = 921.125
$ws.Range("J138").Value = 5268.7085
$ws.Range("K138").Value = 2763.375
$ws.Range("L138").Value = 15806.1255
$ws.Range("M138").Value = 2376.625
$ws.Range("N138").Value = -26086.1255
$ws.Range("H141").Value = 1839.0476
$ws.Range("I141").Value = 1300.3667
$ws.Range("J141").Value = 3185.75
$ws.Range("K141").Value = 3901.1001
$ws.Range("L141").Value = 9557.25
$ws.Range("M141").Value = 1278.8999
$ws.Range("N141").Value = -19917.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5071.3477
$ws.Range("I32").Value = 3500.8228
$ws.Range("J32").Value = 14615.308
$ws.Range("K32").Value = 3500.8228
$ws.Range("L32").Value = 14615.308
$ws.Range("M32").Value = -3213.8228
$ws.Range("N32").Value = -15189.308
$ws.Range("H61").Value = 3047.0188
$ws.Range("I61").Value = 3239.5112
$ws.Range("J61").Value = 1964.25
$ws.Range("K61").Value = 3239.5112
$ws.Range("L61").Value = 1964.25
$ws.Range("M61").Value = -3027.5112
$ws.Range("N61").Value = -2388.25
$ws.Range("H74").Value = 1489.25
$ws.Range("I74").Value = 1469.6875
$ws.Range("J74").Value = 1515.3334
$ws.Range("K74").Value = 1469.6875
$ws.Range("L74").Value = 1515.3334
$ws.Range("M74").Value = -595.6875
$ws.Range("N74").Value = -3263.3334
$ws.Range("H77").Value = 1489.25
$ws.Range("I77").Value = 1469.6875
$ws.Range("J77").Value = 1515.3334
$ws.Range("K77").Value = 7348.4375
$ws.Range("L77").Value = 7576.666999999999
$ws.Range("M77").Value = -2980.4375
$ws.Range("N77").Value = -16312.667
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 1758.3334
$ws.Range("I110").Value = 1175
$ws.Range("K110").Value = 1175
$ws.Range("M110").Value = 870
$ws.Range("H136").Value = 3047.0188
$ws.Range("I136").Value = 3239.5112
$ws.Range("J136").Value = 1964.25
$ws.Range("K136").Value = 9718.533599999999
$ws.Range("L136").Value = 5892.75
$ws.Range("M136").Value = -7168.533599999999
$ws.Range("N136").Value = -10992.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13078.571
$ws.Range("I20").Value = 996
$ws.Range("J20").Value = 32712.75
$ws.Range("K20").Value = 996
$ws.Range("L20").Value = 32712.75
$ws.Range("M20").Value = -749
$ws.Range("N20").Value = -33206.75
$ws.Range("H107").Value = 200001900
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 9281.9375
$ws.Range("I134").Value = 14167.889
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 42503.667
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -39968.667
$ws.Range("N134").Value = -14070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 30978
$ws.Range("I10").Value = 11301.333
$ws.Range("J10").Value = 90008
$ws.Range("K10").Value = 11301.333
$ws.Range("L10").Value = 90008
$ws.Range("M10").Value = -11162.333
$ws.Range("N10").Value = -90286
$ws.Range("H31").Value = 231793.53
$ws.Range("I31").Value = 1760.1666
$ws.Range("J31").Value = 603385.9
$ws.Range("K31").Value = 1760.1666
$ws.Range("L31").Value = 603385.9
$ws.Range("M31").Value = -1465.1666
$ws.Range("N31").Value = -603975.9
$ws.Range("H34").Value = 231793.53
$ws.Range("I34").Value = 1760.1666
$ws.Range("J34").Value = 603385.9
$ws.Range("K34").Value = 1760.1666
$ws.Range("L34").Value = 603385.9
$ws.Range("M34").Value = -1558.1666
$ws.Range("N34").Value = -603789.9
$ws.Range("H107").Value = 19608732
$ws.Range("J107").Value = 1386.8572
$ws.Range("L107").Value = 1386.8572
$ws.Range("N107").Value = -5226.8572
$ws.Range("H134").Value = 1940.0667
$ws.Range("I134").Value = 2173.575
$ws.Range("K134").Value = 6520.724999999999
$ws.Range("M134").Value = -3985.724999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 102231.44
$ws.Range("J5").Value = 301216.7
$ws.Range("L5").Value = 903650.1000000001
$ws.Range("N5").Value = -903874.1000000001
$ws.Range("H134").Value = 8755.795
$ws.Range("I134").Value = 8198.5
$ws.Range("J134").Value = 9143.478999999999
$ws.Range("K134").Value = 24595.5
$ws.Range("L134").Value = 27430.437
$ws.Range("M134").Value = -19525.5
$ws.Range("N134").Value = -37570.437
$ws.Range("H135").Value = 102231.44
$ws.Range("J135").Value = 301216.7
$ws.Range("L135").Value = 2710950.3
$ws.Range("N135").Value = -2716020.3
$ws.Range("H136").Value = 13093.637
$ws.Range("I136").Value = 100030
$ws.Range("K136").Value = 300090
$ws.Range("M136").Value = -294990
$ws.Range("H137").Value = 12365041
$ws.Range("I137").Value = 8594.666999999999
$ws.Range("J137").Value = 27810598
$ws.Range("K137").Value = 25784.001
$ws.Range("L137").Value = 83431794
$ws.Range("M137").Value = -20684.001
$ws.Range("N137").Value = -83441994
$ws.Range("H138").Value = 8727.950000000001
$ws.Range("I138").Value = 10627.066
$ws.Range("J138").Value = 3030.6
$ws.Range("K138").Value = 31881.198
$ws.Range("L138").Value = 9091.799999999999
$ws.Range("M138").Value = -26741.198
$ws.Range("N138").Value = -19371.8
$ws.Range("H139").Value = 5566.0312
$ws.Range("I139").Value = 9547.5
$ws.Range("J139").Value = 3177.15
$ws.Range("K139").Value = 28642.5
$ws.Range("L139").Value = 9531.450000000001
$ws.Range("M139").Value = -23502.5
$ws.Range("N139").Value = -19811.45
$ws.Range("H140").Value = 3020
$ws.Range("I140").Value = 3020
$ws.Range("K140").Value = 9060
$ws.Range("M140").Value = -3880
$ws.Range("H141").Value = 22871.834
$ws.Range("I141").Value = 36243.668
$ws.Range("J141").Value = 9500
$ws.Range("K141").Value = 108731.004
$ws.Range("L141").Value = 28500
$ws.Range("M141").Value = -103551.004
$ws.Range("N141").Value = -38860

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 13146.223
$ws.Range("I12").Value = 14800
$ws.Range("J12").Value = 12939.5
$ws.Range("K12").Value = 14800
$ws.Range("L12").Value = 12939.5
$ws.Range("M12").Value = -14660
$ws.Range("N12").Value = -13219.5
$ws.Range("H80").Value = 2772
$ws.Range("I80").Value = 2277.7778
$ws.Range("J80").Value = 3266.2222
$ws.Range("K80").Value = 2277.7778
$ws.Range("L80").Value = 3266.2222
$ws.Range("M80").Value = -1279.7778
$ws.Range("N80").Value = -5262.2222
$ws.Range("H83").Value = 2772
$ws.Range("I83").Value = 2277.7778
$ws.Range("J83").Value = 3266.2222
$ws.Range("K83").Value = 11388.889
$ws.Range("L83").Value = 16331.111
$ws.Range("M83").Value = -6396.888999999999
$ws.Range("N83").Value = -26315.111
$ws.Range("H102").Value = 1311.1143
$ws.Range("I102").Value = 996.36
$ws.Range("J102").Value = 2098
$ws.Range("K102").Value = 996.36
$ws.Range("L102").Value = 2098
$ws.Range("M102").Value = 625.64
$ws.Range("N102").Value = -5342
$ws.Range("H132").Value = 2259.818
$ws.Range("I132").Value = 2114.8823
$ws.Range("J132").Value = 2413.8125
$ws.Range("K132").Value = 6344.646900000001
$ws.Range("L132").Value = 7241.4375
$ws.Range("M132").Value = -3814.646900000001
$ws.Range("N132").Value = -12301.4375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 43003.6
$ws.Range("I19").Value = 2503
$ws.Range("J19").Value = 70004
$ws.Range("K19").Value = 2503
$ws.Range("L19").Value = 70004
$ws.Range("M19").Value = -2333
$ws.Range("N19").Value = -70344
$ws.Range("H61").Value = 1769.2
$ws.Range("J61").Value = 1702.5
$ws.Range("L61").Value = 1702.5
$ws.Range("N61").Value = -2106.5
$ws.Range("H113").Value = 1769.2
$ws.Range("J113").Value = 1702.5
$ws.Range("L113").Value = 1702.5
$ws.Range("N113").Value = -6042.5
$ws.Range("H132").Value = 11619784
$ws.Range("I132").Value = 16196417
$ws.Range("J132").Value = 2177.6155
$ws.Range("K132").Value = 48589251
$ws.Range("L132").Value = 6532.8465
$ws.Range("M132").Value = -48586721
$ws.Range("N132").Value = -11592.8465

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3400.4
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3400.4
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
$ws.Range("H113").Value = 1247.3214
$ws.Range("I113").Value = 1216.9131
$ws.Range("J113").Value = 1387.2
$ws.Range("K113").Value = 3650.7393
$ws.Range("L113").Value = 4161.6
$ws.Range("M113").Value = -1480.7393
$ws.Range("N113").Value = -8501.6
$ws.Range("H122").Value = 1105.8928
$ws.Range("I122").Value = 1102.6
$ws.Range("J122").Value = 1133.3334
$ws.Range("K122").Value = 3307.8
$ws.Range("L122").Value = 3400.0002
$ws.Range("M122").Value = -857.7999999999997
$ws.Range("N122").Value = -8300.0002
$ws.Range("H132").Value = 863.4828
$ws.Range("I132").Value = 609.0476
$ws.Range("J132").Value = 1531.375
$ws.Range("K132").Value = 1827.1428
$ws.Range("L132").Value = 4594.125
$ws.Range("M132").Value = 702.8571999999999
$ws.Range("N132").Value = -9654.125
$ws.Range("H136").Value = 8477046
$ws.Range("I136").Value = 2799.9736
$ws.Range("J136").Value = 23811396
$ws.Range("K136").Value = 8399.9208
$ws.Range("L136").Value = 71434188
$ws.Range("M136").Value = -5849.9208
$ws.Range("N136").Value = -71439288
